# Update the lattice multiplication exercise table cells to the newly
# generated set of problems, while preserving each cell's formatting
# (font size 32) which Word retains automatically when we set Range.Text.
$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$nl = [char]11   # Word manual line break (renders as <w:br/>)

$tbl.Cell(1, 1).Range.Text = "26 x 11" + $nl + "  1    1" + $nl + "  ----" + $nl + "2|    |" + $nl + "6|    |"
$tbl.Cell(1, 2).Range.Text = "11 x 31" + $nl + "  3    1" + $nl + "  ----" + $nl + "1|    |" + $nl + "1|    |"
$tbl.Cell(1, 3).Range.Text = "89 x 75" + $nl + "  7    5" + $nl + "  ----" + $nl + "8|    |" + $nl + "9|    |"
$tbl.Cell(2, 1).Range.Text = "63 x 51" + $nl + "  5    1" + $nl + "  ----" + $nl + "6|    |" + $nl + "3|    |"
$tbl.Cell(2, 2).Range.Text = "52 x 77" + $nl + "  7    7" + $nl + "  ----" + $nl + "5|    |" + $nl + "2|    |"
$tbl.Cell(2, 3).Range.Text = "64 x 50" + $nl + "  5    0" + $nl + "  ----" + $nl + "6|    |" + $nl + "4|    |"
$tbl.Cell(3, 1).Range.Text = "49 x 87" + $nl + "  8    7" + $nl + "  ----" + $nl + "4|    |" + $nl + "9|    |"
$tbl.Cell(3, 2).Range.Text = "67 x 15" + $nl + "  1    5" + $nl + "  ----" + $nl + "6|    |" + $nl + "7|    |"
$tbl.Cell(3, 3).Range.Text = "19 x 84" + $nl + "  8    4" + $nl + "  ----" + $nl + "1|    |" + $nl + "9|    |"
$tbl.Cell(4, 1).Range.Text = "48 x 84" + $nl + "  8    4" + $nl + "  ----" + $nl + "4|    |" + $nl + "8|    |"
$tbl.Cell(4, 2).Range.Text = "36 x 90" + $nl + "  9    0" + $nl + "  ----" + $nl + "3|    |" + $nl + "6|    |"
$tbl.Cell(4, 3).Range.Text = "49 x 92" + $nl + "  9    2" + $nl + "  ----" + $nl + "4|    |" + $nl + "9|    |"
$tbl.Cell(5, 1).Range.Text = "25 x 93" + $nl + "  9    3" + $nl + "  ----" + $nl + "2|    |" + $nl + "5|    |"
$tbl.Cell(5, 2).Range.Text = "12 x 34" + $nl + "  3    4" + $nl + "  ----" + $nl + "1|    |" + $nl + "2|    |"
$tbl.Cell(5, 3).Range.Text = "82 x 84" + $nl + "  8    4" + $nl + "  ----" + $nl + "8|    |" + $nl + "2|    |"
